# Refresh the "Price" (D) and "Volume(1h)" (E) columns on the cryptos sheet
# with the latest scraped snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.902.65'
$ws.Range("E2").Value = '  -2.19%  '

$ws.Range("D3").Value = '1.791.97'
$ws.Range("E3").Value = '  -1.90%  '

$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.73'
$ws.Range("E5").Value = '  -1.95%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5125'
$ws.Range("E7").Value = '  -0.74%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3933'
$ws.Range("E8").Value = '  +1.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07807'
$ws.Range("E9").Value = '  -7.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.092'
$ws.Range("E10").Value = '  -2.64%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '40.79'
$ws.Range("E11").Value = '  -2.81%  '

$ws.Range("E12").Value = '  -2.85%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.001'
$ws.Range("E13").Value = '  -0.24%  '

$ws.Range("E14").Value = '  -5.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.238'
$ws.Range("E15").Value = '  -4.00%  '

$ws.Range("D16").Value = '1.781.46'
$ws.Range("E16").Value = '  -2.46%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.65'
$ws.Range("E17").Value = '  -2.95%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001078'
$ws.Range("E18").Value = '  -4.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06530'
$ws.Range("E19").Value = '  -1.38%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.0000'
$ws.Range("E20").Value = '  -0.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.09'
$ws.Range("E21").Value = '  -3.95%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.930'
$ws.Range("E22").Value = '  -2.59%  '

$ws.Range("D23").Value = '27.977.42'
$ws.Range("E23").Value = '  -2.07%  '

$ws.Range("E24").Value = '  -3.55%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.230'
$ws.Range("E25").Value = '  -1.91%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.63'
$ws.Range("E26").Value = '  +0.39%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.34'
$ws.Range("E27").Value = '  -4.27%  '

$ws.Range("D28").Value = '1.992.93'
$ws.Range("E28").Value = '  -2.19%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.369'
$ws.Range("E29").Value = '  -1.66%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.73'
$ws.Range("E30").Value = '  +1.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1080'
$ws.Range("E31").Value = '  -1.52%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.043'
$ws.Range("E32").Value = '  -5.15%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.625'
$ws.Range("E33").Value = '  -1.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.491'
$ws.Range("E34").Value = '  -4.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07073'
$ws.Range("E35").Value = '  -8.62%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.873'
$ws.Range("E36").Value = '  +1.41%  '

$ws.Range("E37").Value = '  -3.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2129'
$ws.Range("E38").Value = '  -4.53%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.52'
$ws.Range("E39").Value = '  +0.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.014'
$ws.Range("E40").Value = '  -4.71%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6100'
$ws.Range("E41").Value = '  -4.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9998'
$ws.Range("E42").Value = '  -0.25%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.153'
$ws.Range("E43").Value = '  -3.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.14'
$ws.Range("E44").Value = '  -3.27%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.311'
$ws.Range("E45").Value = '  -6.38%  '

$ws.Range("E46").Value = '  -2.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.714'
$ws.Range("E47").Value = '  -1.73%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.67'
$ws.Range("E48").Value = '  -2.15%  '

$ws.Range("E49").Value = '  +0.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.913'
$ws.Range("E50").Value = '  -4.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06816'
$ws.Range("E51").Value = '  -2.41%  '
